# Weekly update: insert a new data row (row 70) into the "Poroto granado"
# price sheet, shifting the existing rows 70-99 down to 71-100, and fill
# the new row with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70 (default shift = shift cells down),
# which also pushes the sheet's used range / dimension from R99 to R100.
$ws.Rows("70").Insert()

# Populate the newly inserted row 70 with this week's record.
$ws.Range("A70").Value = 2
$ws.Range("B70").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44917
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100112030
$ws.Range("G70").Value = "Poroto granado"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 700
$ws.Range("K70").Value = 16000
$ws.Range("L70").Value = 17000
$ws.Range("M70").Value = 16500
$ws.Range("N70").Value = "$/caja 15 kilos"
$ws.Range("O70").Value = "Provincia de Limarí"
$ws.Range("P70").Value = 1100
$ws.Range("Q70").Value = 15
$ws.Range("R70").Value = "Hortaliza"
